# Populate the "surat jalan" / invoice template with dynamic data.
#
# NOTE on this runtime's Find.Execute: Range.Find.Execute searches
# forward from the Range's Start position all the way to the physical
# end of the document -- it does NOT stop at the Range's End, even with
# Wrap = wdFindStop (0). Because several of the old values here (single
# digits, repeated "Rp ..." amounts, a repeated company name/phone
# number) occur more than once in the document, a Find/Replace scoped
# to "just this paragraph" is not actually safe.
#
# Instead we replace text surgically: for each target paragraph we take
# its Range, trim off the trailing control characters Word reports as
# part of a paragraph's Range.Text (a paragraph-mark \r, and for the
# last paragraph in a table cell also a cell-mark \x07), optionally
# skip a known, constant label prefix (e.g. "Phone: "), and then assign
# .Text on that exact sub-range. That performs an in-place replacement
# of just those characters (preserving the run formatting, e.g. bold)
# without invoking Find's document-wide scan.

$d = $word.ActiveDocument
$nbsp = [char]0x00A0

function Get-VisibleLength($text) {
    $len = $text.Length
    while ($len -gt 0) {
        $code = [int][char]$text[$len - 1]
        if ($code -eq 13 -or $code -eq 7) {
            $len = $len - 1
        } else {
            break
        }
    }
    return $len
}

# Replace the paragraph's visible text (optionally skipping a fixed
# label prefix that must be left untouched) with $new.
function Set-ParagraphValue($index, $prefix, $new) {
    $p = $d.Paragraphs.Item($index)
    $r = $p.Range
    $text = $r.Text
    $visibleLen = Get-VisibleLength $text
    $prefixLen = $prefix.Length
    $start = $r.Start + $prefixLen
    $end = $r.Start + $visibleLen
    $target = $d.Range($start, $end)
    $target.Text = $new
}

# --- Header (paragraphs 6-7) ---
Set-ParagraphValue "6" "INVOICE #" "6"
Set-ParagraphValue "7" "DATE: " "May 8, 2024"

# --- RECIPIENT block (paragraphs 12-15) ---
Set-ParagraphValue "12" "" "Brody Jensen"
Set-ParagraphValue "13" "" "Vaughan Moreno Llc"
Set-ParagraphValue "14" "" "Aut sit in enim et "
Set-ParagraphValue "15" "Phone: " "+1 944 561 4578"

# --- FROM block (paragraphs 17-20) ---
Set-ParagraphValue "17" "" "Admin"
Set-ParagraphValue "18" "" "DEV"
Set-ParagraphValue "19" "" "dev"
Set-ParagraphValue "20" "Phone: " "dev"

# --- Line item 1 (paragraphs 33-37) ---
Set-ParagraphValue "33" "" "3"
Set-ParagraphValue "34" "" "Testing Sample 3"
Set-ParagraphValue "36" "" ("Rp" + $nbsp + "20.000")
Set-ParagraphValue "37" "" ("Rp" + $nbsp + "60.000")

# --- Line item 2 (paragraphs 39-43) ---
Set-ParagraphValue "39" "" "4"
Set-ParagraphValue "40" "" "Testing Sample 2"
Set-ParagraphValue "42" "" ("Rp" + $nbsp + "30.000")
Set-ParagraphValue "43" "" ("Rp" + $nbsp + "120.000")

# --- Line item 3 (paragraphs 45-49) ---
Set-ParagraphValue "45" "" "4"
Set-ParagraphValue "46" "" "Testing Sample Category"
Set-ParagraphValue "48" "" ("Rp" + $nbsp + "40.000")
Set-ParagraphValue "49" "" ("Rp" + $nbsp + "160.000")

# --- Totals table (paragraphs 84, 96) ---
Set-ParagraphValue "84" "" ("Rp" + $nbsp + "340.000")
Set-ParagraphValue "96" "" ("Rp" + $nbsp + "377.400")

# --- Footer (paragraphs 99-100) ---
Set-ParagraphValue "99" "Make all checks payable to " "DEV"
Set-ParagraphValue "100" "If you have any questions concerning this invoice, contact " "dev"
